$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the raw benchmark sample data in column C (new test run numbers) ---
# Test Run 1 (rows 2-5)
$ws.Range("C2").Value = 1795928
$ws.Range("C3").Value = 974073
$ws.Range("C4").Value = 1035084
$ws.Range("C5").Value = 1626093

# Test Run 2 (rows 7-10)
$ws.Range("C7").Value = 652631
$ws.Range("C8").Value = 767637
$ws.Range("C9").Value = 1193956
$ws.Range("C10").Value = 1645465

# Test Run 3 (rows 12-15)
$ws.Range("C12").Value = 378355
$ws.Range("C13").Value = 725574
$ws.Range("C14").Value = 1005584
$ws.Range("C15").Value = 1666428

# Test Run 4 (rows 17-20)
$ws.Range("C17").Value = 557455
$ws.Range("C18").Value = 804092
$ws.Range("C19").Value = 1286410
$ws.Range("C20").Value = 1665436

# Test Run 5 (rows 22-25)
$ws.Range("C22").Value = 538151
$ws.Range("C23").Value = 778896
$ws.Range("C24").Value = 1233067
$ws.Range("C25").Value = 1387094
